$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '29.773.83'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.20%  '

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.889.48'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.74%  '

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.23%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7523'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +3.36%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '239.47'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -1.28%  '

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.19%  '

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.887.82'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.24%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3037'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -2.45%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.40'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -3.04%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06808'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.95%  '

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07931'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.10%  '

$ws.Cells.Item(13, 2).Value = 'Polygon'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7424'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -4.09%  '

$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.882.39'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -1.12%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.141'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -2.24%  '

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '90.36'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.54%  '

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '29.779.11'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -1.28%  '

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.86'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.64%  '

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.907'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.58%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '241.80'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.75%  '

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000007663'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.90%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.13%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.24%  '

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.891'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -2.19%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '165.63'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.58%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.184'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.69%  '

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.60'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -1.95%  '

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1274'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.36%  '

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.013'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -1.49%  '

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.381'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.25%  '

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.511'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.65%  '

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.234'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.24%  '

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.003'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.51%  '

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05193'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.76%  '

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.248'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -2.64%  '

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7237'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -1.62%  '

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.705'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.64%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.01900'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.28%  '

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.765'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.53%  '

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.122'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -3.42%  '

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4376'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.93%  '

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.18'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -4.31%  '

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9999'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.09%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.879'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -2.53%  '

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8249'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -1.01%  '

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.568'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.21%  '

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '99.23'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -1.55%  '

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.679'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.24%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.040.83'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.30%  '

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '35.90'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -4.38%  '

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05940'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.48%  '
